# Apply the "matrix and demand plots" changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Row 3 ("Trips/Day") input values
$ws.Range("B3").Value = 90
$ws.Range("C3").Value = 60
$ws.Range("D3").Value = 120

# Widen columns B:D to fit new currency-formatted numbers
$ws.Range("B:D").EntireColumn.AutoFit()

# Apply a Currency number format (0 decimals) to the aggregate cost rows
$ws.Range("B27:D28").Style = "Currency"
$ws.Range("B27:D28").NumberFormat = "_(""$""* #,##0_);_(""$""* \(#,##0\);_(""$""* ""-""??_);_(@_)"
